$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the six trailing bullet paragraphs (immigration, energy,
#    subsidies, dioxins, pfoa, iris/epa) - only the first paragraph
#    (about the two literature strands) survives.
# ------------------------------------------------------------------
if ($d.Paragraphs.Count -gt 1) {
    $startDel = $d.Paragraphs(2).Range.Start
    $endDel = $d.Paragraphs($d.Paragraphs.Count).Range.End
    $delRange = $d.Range($startDel, $endDel)
    $delRange.Delete()
}

# ------------------------------------------------------------------
# 2. Blank out the two runs of the remaining paragraph, turning them
#    into single spaces (do this before the hyperlink is inserted so
#    the character offsets are not affected by hidden field codes).
# ------------------------------------------------------------------
$firstRunText = "You need to review the t"
$secondRunText = "wo strands of the literature which are in the folder named social, in particular the papers highlighted in green in the word file called summaries. "

$r1 = $d.Range(0, $firstRunText.Length)
$r1.Text = " "

$r2 = $d.Range(1, 1 + $secondRunText.Length)
$r2.Text = " "

# ------------------------------------------------------------------
# 3. Insert the Google Books hyperlink at the very start of the
#    paragraph; its visible text is the URL itself.
# ------------------------------------------------------------------
$url = "https://books.google.dk/books?id=Sp30OtqPCe0C&printsec=frontcover&dq=inauthor:%22Dorthe+Bjerrum+Jensen%22&hl=da&sa=X&redir_esc=y#v=onepage&q&f=false"
$insertionPoint = $d.Range(0, 0)
$d.Hyperlinks.Add($insertionPoint, $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $url) | Out-Null
